$d = $word.ActiveDocument

$pairs = @(
    @("50×16=", "17×32="),
    @("78×18=", "32×69="),
    @("41×20=", "87×60="),
    @("39×81=", "14×99="),
    @("45×55=", "78×59="),
    @("85×15=", "83×92="),
    @("73×37=", "18×42="),
    @("50×19=", "17×72="),
    @("91×43=", "35×44="),
    @("90×26=", "81×39="),
    @("37×29=", "67×86="),
    @("33×29=", "14×49="),
    @("14×88=", "42×28="),
    @("67×95=", "38×14="),
    @("76×94=", "76×83="),
    @("13×85=", "60×49="),
    @("70×59=", "39×74="),
    @("51×56=", "71×90="),
    @("98×23=", "71×12="),
    @("58×54=", "31×23="),
    @("64×67=", "55×73="),
    @("14×32=", "33×13="),
    @("80×49=", "13×85="),
    @("56×49=", "56×85="),
    @("99×14=", "95×88=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
